$d = $word.ActiveDocument

# 1) Change the text of the first paragraph.
$newFirstParaText = "New changes."
$d.Content.Find.Execute("Abfgwbsdbqh e", $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newFirstParaText, 2)

# 2) Remove the paragraphs between the (now retitled) first paragraph and the
#    final "jaidjid. This is changed file" paragraph: "Ioihohqefd", "Aquwhdijqwd",
#    the empty paragraph, "Dqw", "D", "Qwd", "W", "d" (paragraphs 2-9).
#    Deleting the clean paragraph-aligned range first (not touching paragraph 1's
#    own end-of-paragraph mark) keeps the operation reliable.
$pStart = $d.Paragraphs.Item(2)
$pEnd = $d.Paragraphs.Item(9)
$midRange = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$midRange.Delete()

# 3) Merge the first paragraph with what is now the second paragraph
#    ("jaidjid" + ". This is changed file") by deleting the first paragraph's
#    trailing paragraph mark.
$p1 = $d.Paragraphs.Item(1)
$markRange = $d.Range($p1.Range.End - 1, $p1.Range.End)
$markRange.Delete()

# 4) Relocate the "_GoBack" bookmark so it sits right after "New changes."
#    and before "jaidjid" (it previously sat at the very end of the text).
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()
$p1Final = $d.Paragraphs.Item(1)
$bookmarkPos = $p1Final.Range.Start + $newFirstParaText.Length
$newBookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $newBookmarkRange)
